# Apply "Hjemme passive tweaks lichtwark deleted values" edits
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values
$ws.Range("B2").Value = 11.832707150434354
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 12.847658674271434
$ws.Range("E2").Value = 11.843554243009178

# Row 3 values
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 9.9029072152466568
$ws.Range("D3").Value = 11.443310135848185
$ws.Range("E3").Value = 11.414492161629271

# Update the selection on the sheet to match the new edited range
[void]$ws.Range("B1:E3").Select()
